$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a "Role" column (E) with a value per sub-constructor row ---
$ws.Range("E1").Value = "Role"
$ws.Range("E2").Value = "company coordinator"
$ws.Range("E3").Value = "company as"
$ws.Range("E4").Value = "company viewer"

# Give the new column enough room to show the longest role label.
# (ColumnWidth is in characters and gets rounded to the nearest whole
# pixel on write, so this is the closest the engine can land to the
# 30.6640625-character stored width.)
$ws.Columns("E").ColumnWidth = 29.830729166666668

# Restrict the role cells to a dropdown list of the allowed roles,
# mirroring how the login/company-account upload now validates roles.
$ws.Range("E2:E4").Validation.Add(3, 1, 1, '"company coordinator, company as, company viewer"')

# Leave the selection where the author last clicked while testing the upload.
[void]$ws.Range("I11").Select()
